# Auto-generated: rewrite rows 2-10 (sending/target cluster reassignment + new TPM-based metrics)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @("ECs", "Efna5", "Ephb1", "ECs", 2, 0.6666666666666666, 0.114918, 0.344754, 0.04640425382421802, 0.04640425382421801, 3, 1, 1.397474666666667, 4.192424, 0.9253320219903044, 0.9253320219903044, 0.160594993744, 1.445354943696, 0.04293934202011498, 0.04293934202011497)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(2, $c).Value = $row2[$c-1] }

$row3 = @("ECs", "Efna5", "Ephb1", "MuSCs", 2, 0.6666666666666666, 0.114918, 0.344754, 0.04640425382421802, 0.04640425382421801, 3, 1, 0.07803566666666666, 0.234107, 0.05167099121464913, 0.05167099121464913, 0.008967702741999999, 0.08070932467799999, 0.002397753791673518, 0.002397753791673517)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(3, $c).Value = $row3[$c-1] }

$row4 = @("ECs", "Efna5", "Ephb1", "Neutrophils", 2, 0.6666666666666666, 0.114918, 0.344754, 0.04640425382421802, 0.04640425382421801, 1, 0.3333333333333333, 0.034731, 0.104193, 0.02299698679504644, 0.02299698679504644, 0.003991217058, 0.035920953522, 0.001067158012429525, 0.001067158012429525)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(4, $c).Value = $row4[$c-1] }

$row5 = @("FAPs", "Efna5", "Ephb1", "ECs", 3, 1, 2.030023666666667, 6.090071, 0.819730011809897, 0.819730011809897, 3, 1, 1.397474666666667, 4.192424, 0.9253320219903044, 0.9253320219903044, 2.836906646900444, 25.532159822104, 0.7585224293141881, 0.7585224293141881)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(5, $c).Value = $row5[$c-1] }

$row6 = @("FAPs", "Efna5", "Ephb1", "MuSCs", 3, 1, 2.030023666666667, 6.090071, 0.819730011809897, 0.819730011809897, 3, 1, 0.07803566666666666, 0.234107, 0.05167099121464913, 0.05167099121464913, 0.1584142501774444, 1.425728251597, 0.04235626223861342, 0.04235626223861342)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(6, $c).Value = $row6[$c-1] }

$row7 = @("FAPs", "Efna5", "Ephb1", "Neutrophils", 3, 1, 2.030023666666667, 6.090071, 0.819730011809897, 0.819730011809897, 1, 0.3333333333333333, 0.034731, 0.104193, 0.02299698679504644, 0.02299698679504644, 0.070504751967, 0.634542767703, 0.01885132025709547, 0.01885132025709546)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(7, $c).Value = $row7[$c-1] }

$row8 = @("MuSCs", "Efna5", "Ephb1", "ECs", 3, 1, 0.3315123333333334, 0.994537, 0.133865734365885, 0.133865734365885, 3, 1, 1.397474666666667, 4.192424, 0.9253320219903044, 0.9253320219903044, 0.4632800875208889, 4.169520787688, 0.1238702506560014, 0.1238702506560013)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(8, $c).Value = $row8[$c-1] }

$row9 = @("MuSCs", "Efna5", "Ephb1", "MuSCs", 3, 1, 0.3315123333333334, 0.994537, 0.133865734365885, 0.133865734365885, 3, 1, 0.07803566666666666, 0.234107, 0.05167099121464913, 0.05167099121464913, 0.02586978593988889, 0.232828073459, 0.006916975184362198, 0.006916975184362197)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(9, $c).Value = $row9[$c-1] }

$row10 = @("MuSCs", "Efna5", "Ephb1", "Neutrophils", 3, 1, 0.3315123333333334, 0.994537, 0.133865734365885, 0.133865734365885, 1, 0.3333333333333333, 0.034731, 0.104193, 0.02299698679504644, 0.02299698679504644, 0.011513754849, 0.103623793641, 0.003078508525521452, 0.003078508525521451)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(10, $c).Value = $row10[$c-1] }
